$d = $word.ActiveDocument

# Locate the paragraph that ends with "...primary key." so we can insert
# a brand-new paragraph right after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*There is no weak entity since every entity has an id attribute acting as the primary key.*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph"
}

# Insert a new paragraph right after the target paragraph, matching its
# run/paragraph formatting (sz=28 / szCs=28 -> 14pt font size).
$newRange = $target.Range.InsertParagraphAfter()

$insertedPara = $target.Next()
$insertedPara.Range.Text = "All relationships are strong as no entity is weak."
$insertedPara.Range.Font.Size = 14
